$d = $word.ActiveDocument

# 1) "(vbs3) ," -> "(vbs3)," (drop stray space before the comma)
$d.Content.Find.Execute(
    "(vbs3) ,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(vbs3),", 2
)

# 2) typo fix "wstepne" -> "wstępne"
$d.Content.Find.Execute(
    "wstepne",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "wstępne", 2
)

# 3) "mojej pracy" -> "wykonywanej pracy"
$d.Content.Find.Execute(
    "mojej pracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "wykonywanej pracy", 2
)

# 4) "pluton. Który posiada" -> "pluton, który posiada"
$d.Content.Find.Execute(
    "pluton. Który posiada",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pluton, który posiada", 2
)

# 5) Add a new trailing paragraph after the last (SKMB) paragraph, matching its formatting.
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Opisane wyżej moduły plus nieopisany moduł instruktora są łączone w 1 symulator który jest w pełni zsynchronizowany. Dodatkowo symulator umożliwia połączenie kilku stanowisk w grze sieciowej."

Write-Output "edits applied"
